# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# OFF sheet (row 3 = "R" row) - update Short Att, Short Comp, Deep Att, Deep Int
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 203
$wsOff.Range("C3").Value = 143
$wsOff.Range("D3").Value = 50
$wsOff.Range("G3").Value = 3

# DEF sheet (row 3 = "R" row) - update Short Att, Short Comp, Deep Att, Deep Comp
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 181
$wsDef.Range("C3").Value = 135
$wsDef.Range("D3").Value = 36
$wsDef.Range("E3").Value = 15
